# Auto-generated edit script applying odds updates for rows 2, 4, 5, 6
# per the commit diff ("Atualizando o arquivo XLSX").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("Q2").Value = 2.1
$ws.Range("R2").Value = 1.7

# Row 4
$ws.Range("G4").Value = 1.67
$ws.Range("H4").Value = 3.5
$ws.Range("I4").Value = 5.75
$ws.Range("K4").Value = 2.05
$ws.Range("N4").Value = 8
$ws.Range("X4").Value = 7
$ws.Range("AF4").Value = 67
$ws.Range("AI4").Value = 26
$ws.Range("AN4").Value = 3.5
$ws.Range("AO4").Value = 9
$ws.Range("AV4").Value = 67

# Row 5
$ws.Range("G5").Value = 1.4
$ws.Range("H5").Value = 3.75
$ws.Range("I5").Value = 9
$ws.Range("J5").Value = 1.9
$ws.Range("K5").Value = 2.15
$ws.Range("L5").Value = 7.8
$ws.Range("N5").Value = 6.75
$ws.Range("O5").Value = 1.32
$ws.Range("P5").Value = 2.82
$ws.Range("Q5").Value = 1.93
$ws.Range("R5").Value = 1.7
$ws.Range("S5").Value = 1.4
$ws.Range("T5").Value = 2.52
$ws.Range("U5").Value = 2.15
$ws.Range("V5").Value = 1.55
$ws.Range("X5").Value = 5.7
$ws.Range("Z5").Value = 8.75
$ws.Range("AA5").Value = 13
$ws.Range("AC5").Value = 8.25
$ws.Range("AD5").Value = 7.8
$ws.Range("AE5").Value = 23
$ws.Range("AF5").Value = 150
$ws.Range("AH5").Value = 19.5
$ws.Range("AI5").Value = 65
$ws.Range("AJ5").Value = 28
$ws.Range("AK5").Value = 300
$ws.Range("AL5").Value = 150
$ws.Range("AM5").Value = 110
$ws.Range("AN5").Value = 3.05
$ws.Range("AO5").Value = 6.4
$ws.Range("AP5").Value = 18
$ws.Range("AR5").Value = 55
$ws.Range("AS5").Value = 300
$ws.Range("AT5").Value = 2.5
$ws.Range("AU5").Value = 8.5
$ws.Range("AV5").Value = 90
$ws.Range("AW5").Value = 9.25
$ws.Range("AX5").Value = 55
$ws.Range("AY5").Value = 50
$ws.Range("AZ5").Value = 500
$ws.Range("BA5").Value = 450

# Row 6
$ws.Range("G6").Value = 1.75
$ws.Range("H6").Value = 3.75
$ws.Range("I6").Value = 4
$ws.Range("J6").Value = 2.3
$ws.Range("K6").Value = 2.22
$ws.Range("L6").Value = 4.3
$ws.Range("P6").Value = 3.4
$ws.Range("U6").Value = 1.65
$ws.Range("V6").Value = 1.98
$ws.Range("W6").Value = 8
$ws.Range("X6").Value = 9
$ws.Range("Y6").Value = 8.25
$ws.Range("Z6").Value = 14.5
$ws.Range("AA6").Value = 13
$ws.Range("AB6").Value = 23
$ws.Range("AD6").Value = 7.3
$ws.Range("AE6").Value = 14.5
$ws.Range("AF6").Value = 60
$ws.Range("AG6").Value = 400
$ws.Range("AH6").Value = 13
$ws.Range("AI6").Value = 23
$ws.Range("AJ6").Value = 13.5
$ws.Range("AK6").Value = 60
$ws.Range("AL6").Value = 35
$ws.Range("AM6").Value = 37
$ws.Range("AN6").Value = 3.7
$ws.Range("AO6").Value = 8.5
$ws.Range("AP6").Value = 16.5
$ws.Range("AQ6").Value = 28
$ws.Range("AR6").Value = 55
$ws.Range("AU6").Value = 7.1
$ws.Range("AV6").Value = 60
$ws.Range("AW6").Value = 5.8
$ws.Range("AX6").Value = 22
$ws.Range("AY6").Value = 27
$ws.Range("AZ6").Value = 110
$ws.Range("BA6").Value = 150
$ws.Range("BB6").Value = 300
